$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "real" burndown series (row 4): update last tracked day and extend one more day
$ws.Range("U4").Value = 22
$ws.Range("V4").Value = 20

# "finished" row (row 6): update last tracked day and extend one more day
$ws.Range("T6").Value = 6
$ws.Range("U6").Value = 2

# Lilly's week 3 hours (row 14) - formula in F14 recalculates automatically
$ws.Range("E14").Value = 2

# Move the saved selection to I11 (matches author's last-clicked cell)
$ws.Range("I11").Select() | Out-Null
